# Update "想去人数" (number of people interested) figures on gh-pages
# output regeneration (commit 456a3b4).
$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value = 2462
$ws1.Range("F5").Value = 1670
$ws1.Range("F9").Value = 3506
$ws1.Range("F11").Value = 1153
$ws1.Range("F12").Value = 1566
$ws1.Range("F14").Value = 883
$ws1.Range("F15").Value = 13
$ws1.Range("F16").Value = 1260
$ws1.Range("F19").Value = 444
$ws1.Range("F22").Value = 2086
$ws1.Range("F24").Value = 4246

# --- Sheet 2: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F14").Value = 35
$ws2.Range("F15").Value = 35
$ws2.Range("F20").Value = 14
$ws2.Range("F23").Value = 118
$ws2.Range("F44").Value = 83

# --- Sheet 3: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 2538
$ws3.Range("F6").Value = 9572
$ws3.Range("F11").Value = 2956
$ws3.Range("F14").Value = 193

# --- Sheet 4: 全部类型 (All Types, aggregated) ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 2538
$ws4.Range("F5").Value = 2462
$ws4.Range("F8").Value = 2956
$ws4.Range("F17").Value = 1153
$ws4.Range("F20").Value = 883
$ws4.Range("F23").Value = 1260
$ws4.Range("F26").Value = 35
$ws4.Range("F27").Value = 35
$ws4.Range("F29").Value = 14
$ws4.Range("F31").Value = 444
$ws4.Range("F35").Value = 118
$ws4.Range("F36").Value = 118
$ws4.Range("F40").Value = 2086
$ws4.Range("F44").Value = 4246
